# Apply updated loading_percent values (Case_3_50, 380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.29126840243755
$ws.Range("C2").Value = 11.40952840570048
$ws.Range("D2").Value = 4.515686686952149
$ws.Range("F2").Value = 26.68701402281247
$ws.Range("G2").Value = 3.621311576494712
$ws.Range("I2").Value = 20.23164453289785
$ws.Range("N2").Value = 16.45933191183161
$ws.Range("B3").Value = 14.64272053545449
$ws.Range("C3").Value = 10.75739722647943
$ws.Range("D3").Value = 4.539222349931217
$ws.Range("F3").Value = 26.4889307128727
$ws.Range("G3").Value = 3.624713716211601
$ws.Range("I3").Value = 20.29201184188097
$ws.Range("N3").Value = 16.53167223986506
$ws.Range("B4").Value = 14.23291915474673
$ws.Range("C4").Value = 10.33837539094934
$ws.Range("D4").Value = 4.554286793787387
$ws.Range("F4").Value = 26.37784715841357
$ws.Range("G4").Value = 3.62691043014828
$ws.Range("I4").Value = 20.33597824671418
$ws.Range("N4").Value = 16.57810081183548
$ws.Range("B5").Value = 14.06327191036689
$ws.Range("C5").Value = 10.16310787339057
$ws.Range("D5").Value = 4.560580578088538
$ws.Range("F5").Value = 26.33526996963107
$ws.Range("G5").Value = 3.627832813183698
$ws.Range("I5").Value = 20.35561745376909
$ws.Range("N5").Value = 16.59752845191983
$ws.Range("B6").Value = 14.03495043987692
$ws.Range("C6").Value = 10.13373783896007
$ws.Range("D6").Value = 4.561635030752857
$ws.Range("F6").Value = 26.32836358468564
$ws.Range("G6").Value = 3.627987620206557
$ws.Range("I6").Value = 20.35898218497959
$ws.Range("N6").Value = 16.60078511040488
$ws.Range("B7").Value = 14.23064158841096
$ws.Range("C7").Value = 10.33602969781891
$ws.Range("D7").Value = 4.554371045950036
$ws.Range("F7").Value = 26.37726200701284
$ws.Range("G7").Value = 3.626922759447032
$ws.Range("I7").Value = 20.33623615071441
$ws.Range("N7").Value = 16.57836076202941
$ws.Range("B8").Value = 15.07020358243731
$ws.Range("C8").Value = 11.18862635578854
$ws.Range("D8").Value = 4.523674850138188
$ws.Range("F8").Value = 26.61655301320004
$ws.Range("G8").Value = 3.622462323833028
$ws.Range("I8").Value = 20.25101949416507
$ws.Range("N8").Value = 16.48385852933076
$ws.Range("B9").Value = 16.61457399712624
$ws.Range("C9").Value = 12.70711200977432
$ws.Range("D9").Value = 4.468317507394975
$ws.Range("F9").Value = 27.16738798069339
$ws.Range("G9").Value = 3.614566070939022
$ws.Range("I9").Value = 20.13920452868437
$ws.Range("N9").Value = 16.31441324386573
$ws.Range("B10").Value = 17.67561893838188
$ws.Range("C10").Value = 13.72321612493207
$ws.Range("D10").Value = 4.4305538007949
$ws.Range("F10").Value = 27.61868825145148
$ws.Range("G10").Value = 3.609276843331576
$ws.Range("I10").Value = 20.09146315121391
$ws.Range("N10").Value = 16.19947756897527
$ws.Range("B11").Value = 18.14044135778222
$ws.Range("C11").Value = 14.16299946188439
$ws.Range("D11").Value = 4.413996483674102
$ws.Range("F11").Value = 27.83334615284807
$ws.Range("G11").Value = 3.606980471960467
$ws.Range("I11").Value = 20.07735147551886
$ws.Range("N11").Value = 16.14923950504755
$ws.Range("B12").Value = 18.31376260010654
$ws.Range("C12").Value = 14.32625438316412
$ws.Range("D12").Value = 4.40781537332994
$ws.Range("F12").Value = 27.91590697784245
$ws.Range("G12").Value = 3.606126568912492
$ws.Range("I12").Value = 20.07311139209957
$ws.Range("N12").Value = 16.13050804776085
$ws.Range("B13").Value = 18.27655657606698
$ws.Range("C13").Value = 14.29124118329708
$ws.Range("D13").Value = 4.409142646666135
$ws.Range("F13").Value = 27.89807050972057
$ws.Range("G13").Value = 3.606309776222473
$ws.Range("I13").Value = 20.07397533293528
$ws.Range("N13").Value = 16.13452921555416
$ws.Range("B14").Value = 18.15475532519791
$ws.Range("C14").Value = 14.17649657163892
$ws.Range("D14").Value = 4.413486183993016
$ws.Range("F14").Value = 27.84011335746405
$ws.Range("G14").Value = 3.606909907077747
$ws.Range("I14").Value = 20.07698046745742
$ws.Range("N14").Value = 16.14769260225962
$ws.Range("B15").Value = 18.07979380545355
$ws.Range("C15").Value = 14.10578345844781
$ws.Range("D15").Value = 4.416158271018085
$ws.Range("F15").Value = 27.8047767407881
$ws.Range("G15").Value = 3.607279544115959
$ws.Range("I15").Value = 20.07896521475335
$ws.Range("N15").Value = 16.1557936122111
$ws.Range("B16").Value = 17.64487091937854
$ws.Range("C16").Value = 13.69401868219077
$ws.Range("D16").Value = 4.431648317103405
$ws.Range("F16").Value = 27.60484216196934
$ws.Range("G16").Value = 3.60942911638578
$ws.Range("I16").Value = 20.09253930655924
$ws.Range("N16").Value = 16.20280178366986
$ws.Range("B17").Value = 17.37338870251203
$ws.Range("C17").Value = 13.43562279655227
$ws.Range("D17").Value = 4.441309734427143
$ws.Range("F17").Value = 27.48453620031837
$ws.Range("G17").Value = 3.610775844528149
$ws.Range("I17").Value = 20.10282202350594
$ws.Range("N17").Value = 16.23216277243153
$ws.Range("B18").Value = 17.21556538499416
$ws.Range("C18").Value = 13.28488976336113
$ws.Range("D18").Value = 4.446925267713351
$ws.Range("F18").Value = 27.41622481714136
$ws.Range("G18").Value = 3.611560780428304
$ws.Range("I18").Value = 20.10945175021748
$ws.Range("N18").Value = 16.24924318799629
$ws.Range("B19").Value = 17.16184590475504
$ws.Range("C19").Value = 13.23349349876926
$ws.Range("D19").Value = 4.448836663626182
$ws.Range("F19").Value = 27.39324994748705
$ws.Range("G19").Value = 3.611828323935461
$ws.Range("I19").Value = 20.11181901207581
$ws.Range("N19").Value = 16.25505948087556
$ws.Range("B20").Value = 17.4024626958358
$ws.Range("C20").Value = 13.463348341362
$ws.Range("D20").Value = 4.440275205625118
$ws.Range("F20").Value = 27.49725184382937
$ws.Range("G20").Value = 3.610631414235891
$ws.Range("I20").Value = 20.10165329716157
$ws.Range("N20").Value = 16.22901730858116
$ws.Range("B21").Value = 18.1906054365268
$ws.Range("C21").Value = 14.21028927025823
$ws.Range("D21").Value = 4.412207976877713
$ws.Range("F21").Value = 27.85710277107706
$ws.Range("G21").Value = 3.606733209224617
$ws.Range("I21").Value = 20.07606775728936
$ws.Range("N21").Value = 16.14381826879696
$ws.Range("B22").Value = 18.68993563732032
$ws.Range("C22").Value = 14.67931469574046
$ws.Range("D22").Value = 4.394381662515722
$ws.Range("F22").Value = 28.09968000028968
$ws.Range("G22").Value = 3.604276875201072
$ws.Range("I22").Value = 20.06578336355483
$ws.Range("N22").Value = 16.08984052581338
$ws.Range("B23").Value = 18.42491430763704
$ws.Range("C23").Value = 14.43075327188616
$ws.Range("D23").Value = 4.403848769755451
$ws.Range("F23").Value = 27.96955975766434
$ws.Range("G23").Value = 3.605579538083882
$ws.Range("I23").Value = 20.07068026117194
$ws.Range("N23").Value = 16.11849403991833
$ws.Range("B24").Value = 17.38932376735417
$ws.Range("C24").Value = 13.45082040180503
$ws.Range("D24").Value = 4.440742725846909
$ws.Range("F24").Value = 27.49150043550759
$ws.Range("G24").Value = 3.610696677884021
$ws.Range("I24").Value = 20.10217944239181
$ws.Range("N24").Value = 16.23043874842781
$ws.Range("B25").Value = 16.20898526075382
$ws.Range("C25").Value = 12.31345148663811
$ws.Range("D25").Value = 4.48277957033947
$ws.Range("F25").Value = 27.00994159702672
$ws.Range("G25").Value = 3.616611814980879
$ws.Range("I25").Value = 20.16345883978886
$ws.Range("N25").Value = 16.35856590356731
